$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new POST test-case rows right after the existing POST row (row 3) ---
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = "POST"
$ws.Range("B4").Value = "/redfish/v1/AccountService/Accounts"
$ws.Range("C4").Value = '{"UserName": "test2user", "Password": "Test1234", "RoleId": "Administrator","redundant-key":"redundant-value"}'

$ws.Range("A5").Value = "POST"
$ws.Range("B5").Value = "/redfish/v1/AccountService/Accounts"
$ws.Range("C5").Value = '{"UserName": "test3user", "password": "Test1234", "RoleId": "Administrator"}'

# --- Append a new DELETE test case (non-existent user) as the final row ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1
$ws.Cells.Item($newRow, 1).Value = "DELETE"
$ws.Cells.Item($newRow, 2).Value = "/redfish/v1/AccountService/Accounts/testwronguser"

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 43.5 - 0.8333333333333333
$ws.Columns.Item(3).ColumnWidth = 88 - 0.8333333333333333

# --- Selection ---
$ws.Range("C18").Select() | Out-Null
